# Conserto do erro com o rótulo da coluna 2050 nas tabelas e retirada das
# linhas com total das tabelas.
#
# For each of the first four tables (and the 5th, which has no Total row)
# the E1 header cell incorrectly held a stray numeric value
# (792.3439386676173) instead of the "2050" (or "2041-2050") text label
# used by the other header cells in that row. We fix the label and, for
# the tables that have one, drop the trailing "Total" row - along with
# the Custo Total sheet's own "Total" row.

$wb = $excel.ActiveWorkbook

# Sheet 1: Potencia Acumulada - SIN (MW)  -> header 2015/2030/2040/2050
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E1").Value = "'2050"
$ws1.Rows.Item(13).Delete()

# Sheet 2: Geracao Periodo Medio (MWMed) -> header 2015/2030/2040/2050
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("E1").Value = "'2050"
$ws2.Rows.Item(13).Delete()

# Sheet 3: Atendimento a Ponta(MW) -> header 2015/2030/2040/2050
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("E1").Value = "'2050"
$ws3.Rows.Item(13).Delete()

# Sheet 4: Potencia Incremental - SIN(MW) -> header 2015/2015-2030/2031-2040/2041-2050
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("E1").Value = "'2041-2050"
$ws4.Rows.Item(13).Delete()

# Sheet 5: Emissoes Totais (MtCO2eq) -> header 2015/2030/2040/2050 (no Total row here)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("E1").Value = "'2050"

# Sheet 6: Custo Total (bilhões de R$) -> drop its own "Total" row
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
